$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    1.1920823830940299,
    1.44483043833123,
    1.60437696642062,
    1.61721518221868,
    1.9384270008325299,
    2.27500687401893,
    2.5263205924989798,
    2.6763465352767399,
    3.28376077921294,
    3.6560473312344,
    3.9185144530215101,
    4.0153152953595601,
    4.8793481283687097,
    5.2926747649666703,
    5.9354034932481499,
    6.9073644628952504,
    6.9447638559519502,
    7.0303857232502702,
    7.1265532702269896,
    7.3122769616025201,
    7.6250191900396098,
    7.8039675252019096,
    7.8752505944591897,
    7.9038054919662004,
    8.3142087476862692,
    8.4865013288392497,
    8.5518727686377698,
    8.9628417632004904,
    8.9996003726823908,
    9.1972764686211299,
    9.6380718577846807,
    10.034444925142701,
    10.485012707054,
    10.7051760366842,
    11.7893730352987,
    12.1857578694085,
    12.624548570155399,
    13.4241568354416,
    13.9860874403455,
    14.4154195872818,
    14.5545893043604,
    14.802780029186501,
    15.605375290458401,
    16.288983179863202,
    16.393157076385201,
    16.4606452540128,
    16.468570729338499,
    16.939434744378602,
    17.137520182095901,
    17.267756447434401,
    17.4068112998396,
    17.432538869536199,
    18.1049314386863,
    18.348661974797,
    18.680562437735301,
    18.876756438582799,
    19.198544637750398,
    19.488146172722299,
    19.504405800634999,
    19.865343827647902
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

$ws.Range("A2:A61").Select()
